$d = $word.ActiveDocument

function Get-GoBackParagraph($doc) {
    $bm = $doc.Bookmarks.Item("_GoBack")
    $bmRange = $bm.Range
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -eq $bmRange.Start) {
            return $p
        }
    }
    return $null
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Insert new paragraph: "The ai crashes all the time now" ---
$bookmarkPara = Get-GoBackParagraph $d
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara = Get-GoBackParagraph $d
$newParaRange = $bookmarkPara.Previous().Range
$xml1 = $pkgHeader + '<w:body><w:p><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ai</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> crashes all the time now</w:t></w:r></w:p></w:body>' + $pkgFooter
$newParaRange.InsertXML($xml1)

# --- Insert new paragraph: "Cant build in the right places" ---
$bookmarkPara = Get-GoBackParagraph $d
$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara = Get-GoBackParagraph $d
$newParaRange2 = $bookmarkPara.Previous().Range
$xml2 = $pkgHeader + '<w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>Cant</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> build in the right places</w:t></w:r></w:p></w:body>' + $pkgFooter
$newParaRange2.InsertXML($xml2)

# --- Insert new run "Cant really do anything" at the start of the bookmark
#     paragraph, before the bookmarkStart/bookmarkEnd it already contains.
#     Plain-text InsertBefore (rather than InsertXML) is used here so the
#     existing _GoBack bookmark at that position is preserved. ---
$bookmarkPara = Get-GoBackParagraph $d
$insertionRange = $bookmarkPara.Range.Duplicate
$insertionRange.Collapse(1)
$insertionRange.InsertBefore("Cant really do anything")
